# Insert a new weekly price record as row 166 in the "Hortaliza, Femacal de
# La Calera - Berenjena" sheet. All existing rows from 166 downward shift
# down by one (to 167..179) and the new row is populated with the latest
# observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 166:178 down to 167:179, creating a blank row 166.
$ws.Rows.Item(166).Insert()

# Populate the new row 166 with the new weekly record.
$ws.Cells.Item(166, 1).Value = 3
$ws.Cells.Item(166, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(166, 3).Value = "Coquimbo"
$ws.Cells.Item(166, 4).Value = 44516
$ws.Cells.Item(166, 5).Value = 5
$ws.Cells.Item(166, 6).Value = 100112001
$ws.Cells.Item(166, 7).Value = "Berenjena"
$ws.Cells.Item(166, 8).Value = "Sin especificar"
$ws.Cells.Item(166, 9).Value = "Primera"
$ws.Cells.Item(166, 10).Value = 75
$ws.Cells.Item(166, 11).Value = 7500
$ws.Cells.Item(166, 12).Value = 8000
$ws.Cells.Item(166, 13).Value = 7767
$ws.Cells.Item(166, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(166, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(166, 16).Value = 129
$ws.Cells.Item(166, 17).Value = 60
$ws.Cells.Item(166, 18).Value = "Hortaliza"
